$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 90 entirely; Excel automatically shifts subsequent rows up.
$ws.Rows.Item(90).Delete()
